$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = " Sep 27"
$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 5
$ws.Range("C4").Value = 2
